# Fixed the back button error
# The DefectType for a few batches was recorded incorrectly; correct them to "Type C".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BATCH-003 (row 4) was "Type A", should be "Type C"
$ws.Range("C4").Value = "Type C"

# BATCH-007 (row 8) was "Type B", should be "Type C"
$ws.Range("C8").Value = "Type C"

# BATCH-037 (row 38) was "Type B", should be "Type C"
$ws.Range("C38").Value = "Type C"

# Restore the view to where the user left off (scrolled down, selection on F36)
$ws.Application.ActiveWindow.ScrollRow = 27
$ws.Range("F36").Select()
